$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.196254014968872
$ws.Range("B1").Value = 2.482172012329102
$ws.Range("C1").Value = 4.134514808654785
$ws.Range("D1").Value = 2.089826583862305
$ws.Range("E1").Value = 1.183972597122192
